# Updated cryptos list with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns scraped from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.615.21"
$ws.Range("E2").Value = "  +4.75%  "
# Row 3
$ws.Range("D3").Value = "2.724.99"
$ws.Range("E3").Value = "  +2.87%  "
# Row 4
$ws.Range("E4").Value = "  +0.10%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.33"
$ws.Range("E5").Value = "  -0.41%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.28"
$ws.Range("E6").Value = "  +5.93%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.21%  "
# Row 8
$ws.Range("E8").Value = "  +1.35%  "
# Row 9
$ws.Range("D9").Value = "2.749.06"
$ws.Range("E9").Value = "  +3.15%  "
# Row 10
$ws.Range("E10").Value = "  +1.69%  "
# Row 11
$ws.Range("E11").Value = "  +5.33%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.163"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.390"
$ws.Range("E13").Value = "  +3.60%  "
# Row 14
$ws.Range("D14").Value = "3.210.23"
$ws.Range("E14").Value = "  +2.88%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.29"
$ws.Range("E15").Value = "  +1.97%  "
# Row 16
$ws.Range("D16").Value = "63.536.83"
$ws.Range("E16").Value = "  +4.65%  "
# Row 17
$ws.Range("E17").Value = "  +6.09%  "
# Row 18
$ws.Range("D18").Value = "2.740.01"
$ws.Range("E18").Value = "  +3.08%  "
# Row 19
$ws.Range("E19").Value = "  +3.26%  "
# Row 20
$ws.Range("E20").Value = "  +2.67%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.33"
$ws.Range("E21").Value = "  +2.67%  "
# Row 22
$ws.Range("E22").Value = "  -0.19%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.61%  "
# Row 24
$ws.Range("E24").Value = "  +0.49%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.87"
$ws.Range("E25").Value = "  +2.99%  "
# Row 26
$ws.Range("E26").Value = "  +4.15%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.52"
$ws.Range("E27").Value = "  +4.74%  "
# Row 28
$ws.Range("E28").Value = "  -0.11%  "
# Row 29
$ws.Range("E29").Value = "  +11.60%  "
# Row 30
$ws.Range("E30").Value = "  -0.90%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.10"
$ws.Range("E31").Value = "  +6.34%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "171.69"
$ws.Range("E32").Value = "  +1.29%  "
# Row 33
$ws.Range("E33").Value = "  +12.91%  "
# Row 34
$ws.Range("E34").Value = "  -0.12%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.47"
$ws.Range("E35").Value = "  +3.67%  "
# Row 36
$ws.Range("E36").Value = "  +7.04%  "
# Row 37
$ws.Range("E37").Value = "  +9.26%  "
# Row 38
$ws.Range("E38").Value = "  +9.70%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  +13.54%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "344.58"
$ws.Range("E40").Value = "  +4.13%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  +5.39%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.24"
$ws.Range("E42").Value = "  +2.17%  "
# Row 43
$ws.Range("E43").Value = "  +6.58%  "
# Row 44
$ws.Range("E44").Value = "  +6.27%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.77"
$ws.Range("E45").Value = "  +3.77%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0591"
$ws.Range("E46").Value = "  +5.30%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "139.58"
$ws.Range("E47").Value = "  +4.19%  "
# Row 48
$ws.Range("E48").Value = "  +4.31%  "
# Row 49
$ws.Range("E49").Value = "  +3.21%  "
# Row 50
$ws.Range("E50").Value = "  +0.91%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.996"
$ws.Range("E51").Value = "  -0.28%  "
